$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Line data")

# Fix power ratings for the two lines from bus 1 to 2 (rows 2 and 3).
$ws.Range("F2").Value = 47.5
$ws.Range("F3").Value = 47.5
